$wb = $excel.ActiveWorkbook

$renames = @(
    @{ Old = "summ32807652"; New = "summ16711027" },
    @{ Old = "summ33005472"; New = "summ17116329" },
    @{ Old = "summ33227580"; New = "summ17353373" },
    @{ Old = "summ33464105"; New = "summ17589307" },
    @{ Old = "summ33703834"; New = "summ17833707" },
    @{ Old = "summ33940790"; New = "summ18069531" },
    @{ Old = "summ34173058"; New = "summ18309973" },
    @{ Old = "summ34409271"; New = "summ18533066" },
    @{ Old = "summ34661093"; New = "summ18782861" }
)

foreach ($pair in $renames) {
    $sheet = $wb.Worksheets.Item($pair.Old)
    $sheet.Name = $pair.New
}
